$wb = $excel.ActiveWorkbook

# --- Sheet "Summary" ---
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B2").Value = 0.6226591760299626
$summary.Range("C2").Value = 0.5771495877502945
$summary.Range("D2").Value = 0.9176029962546817
$summary.Range("E2").Value = 0.7086044830079538
$summary.Range("F2").Value = 0.8207705192629816
$summary.Range("G2").Value = 0.8972462849496443
$summary.Range("H2").Value = 0.7754176661195977
$summary.Range("I2").Value = 490
$summary.Range("J2").Value = 359
$summary.Range("K2").Value = 175
$summary.Range("L2").Value = 44

# --- Sheet "Classification Report" ---
$clf = $wb.Worksheets.Item("Classification Report")
$clf.Range("B2").Value = 0.7990867579908676
$clf.Range("C2").Value = 0.3277153558052435
$clf.Range("D2").Value = 0.4648074369189907

$clf.Range("B3").Value = 0.5771495877502945
$clf.Range("C3").Value = 0.9176029962546817
$clf.Range("D3").Value = 0.7086044830079538

$clf.Range("B4").Value = 0.6226591760299626
$clf.Range("C4").Value = 0.6226591760299626
$clf.Range("D4").Value = 0.6226591760299626
$clf.Range("E4").Value = 0.6226591760299626

$clf.Range("B5").Value = 0.688118172870581
$clf.Range("C5").Value = 0.6226591760299626
$clf.Range("D5").Value = 0.5867059599634722

$clf.Range("B6").Value = 0.6881181728705811
$clf.Range("C6").Value = 0.6226591760299626
$clf.Range("D6").Value = 0.5867059599634722

# --- Sheet "Confusion Matrix" ---
$cm = $wb.Worksheets.Item("Confusion Matrix")
$cm.Range("B2").Value = 175
$cm.Range("C2").Value = 359
$cm.Range("B3").Value = 44
$cm.Range("C3").Value = 490
